# Upgrade database schema migration
# Append a new data row (row 29) to each of the four worksheets, mirroring
# the existing layout used by the prior rows (row 28, etc.).

$wb = $excel.ActiveWorkbook

# Data for the new row, per worksheet name.
# Columns: A(time) B(总长) C(ID) D(实际长度) E(和校验) F(总长_DEC) G(ID_DEC) H(实际长度_DEC) I(和校验_DEC)
$rowsToAdd = @{
    "DE_LFT_#1" = @{
        A = 45815.43923611111
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x78"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 376
        I = 14
    }
    "DE_LFT_#2" = @{
        A = 45815.43923611111
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x78"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 376
        I = 14
    }
    "DE_PLT_#1" = @{
        A = 45815.43923611111
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x82"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 129
        I = 7
    }
    "DE_PLT_#2" = @{
        A = 45815.43923611111
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x81"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 129
        I = 3
    }
}

$newRow = 29

foreach ($sheetName in $rowsToAdd.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rowsToAdd[$sheetName]

    $ws.Range("A$newRow").Value = $data.A
    $ws.Range("A$newRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B$newRow").Value = $data.B
    $ws.Range("C$newRow").Value = $data.C
    $ws.Range("D$newRow").Value = $data.D
    $ws.Range("E$newRow").Value = $data.E

    $ws.Range("F$newRow").Value = $data.F
    $ws.Range("G$newRow").Value = $data.G
    $ws.Range("H$newRow").Value = $data.H
    $ws.Range("I$newRow").Value = $data.I
}
